# Append two new job postings to the "ランサーズ" sheet and refresh the
# "取得日時" (fetched-at) timestamp on every row to 2025-11-28 12:36:47.
#
# Final row layout (1-indexed, row 1 = header):
#   2 (NEW)   爆サイ AI 開発者募集
#   3 (was 2) エクセルデータ転記作業の効率化依頼
#   4 (NEW)   MT5アラートツールの制作
#   5 (was 3) コンサル会社のバックオフィス業務フロー設計 一括見積依頼
#   6 (was 4) PSE認証代行をお手伝いしてくれる方募集!

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Drop the existing hyperlink objects before shifting rows around --
# Insert() does not re-target <hyperlinks> refs, so stale links would
# otherwise stay pinned to their original cells.
$ws.Range("F2").Hyperlinks.Delete()
$ws.Range("F3").Hyperlinks.Delete()
$ws.Range("F4").Hyperlinks.Delete()

# Make room for the two new listings: one above the current row 2, and
# one above what is (after the first insert) row 4.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(4).Insert()

# --- refresh the fetch timestamp on every existing (shifted) row ---
$ws.Range("A3").Value = "2025-11-28 12:36:47"
$ws.Range("A5").Value = "2025-11-28 12:36:47"
$ws.Range("A6").Value = "2025-11-28 12:36:47"

# --- new row 2: 爆サイ AI 開発者募集 ---
$ws.Range("A2").Value = "2025-11-28 12:36:47"
$ws.Range("B2").Value = "【急募】掲示板サイト(爆サイ)でAIによる自然な会話で書き込みを埋めていけるソフト開発者募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5443464"
$ws.Range("G2").Value = 378
$ws.Range("H2").Value = "🔥AI,Ai ◆開発 ◇サイト"

# --- new row 4: MT5アラートツールの制作 ---
$ws.Range("A4").Value = "2025-11-28 12:36:47"
$ws.Range("B4").Value = "MT5アラートツールの制作"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5443470"
$ws.Range("G4").Value = 65
$ws.Range("H4").Value = "◆ツール"

# --- re-create the URL hyperlinks for every data row (F2:F6) ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5443464")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5442971")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5443470")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5442904")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5443188")

# --- widen column H (スキル概要) from 12 to 17 characters ---
# ColumnWidth is in "characters"; the engine stores width with Excel's
# usual +5/6 padding baked in, so request width-5/6 to land exactly on 17.
$ws.Columns.Item(8).ColumnWidth = 16.16666666666667
